# Regenerate save_data to use K (column G) instead of Strike# and update its values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" values keyed by row number (column G), computed/recalculated upstream.
$kValues = @{
    3  = 1
    4  = 2
    5  = 2
    6  = 0
    7  = 2
    8  = 3
    9  = 2
    10 = 1
    11 = 0
    13 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
